$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 70538
$ws.Range("E2").Value = 1386381840795
$ws.Range("F2").Value = 17271056717
$ws.Range("G2").Value = 0.45358

$ws.Range("D3").Value = 3618.96
$ws.Range("E3").Value = 434252004740
$ws.Range("F3").Value = 10797253507
$ws.Range("G3").Value = 1.67228

$ws.Range("D4").Value = 1.001
$ws.Range("E4").Value = 104525426959
$ws.Range("F4").Value = 35083851689
$ws.Range("G4").Value = 0.06415999999999999

$ws.Range("D5").Value = 606.34
$ws.Range("E5").Value = 93283807208
$ws.Range("F5").Value = 1111850257
$ws.Range("G5").Value = 0.18988

$ws.Range("D6").Value = 196.5
$ws.Range("E6").Value = 87275375056
$ws.Range("F6").Value = 2264556372
$ws.Range("G6").Value = -0.70403

$ws.Range("B7").Value = 'STETH'
$ws.Range("C7").Value = 'Lido Staked Ether'
$ws.Range("D7").Value = 3604.15
$ws.Range("E7").Value = 34934850880
$ws.Range("F7").Value = 57655138
$ws.Range("G7").Value = 1.56875

$ws.Range("B8").Value = 'XRP'
$ws.Range("C8").Value = 'XRP'
$ws.Range("D8").Value = 0.627454
$ws.Range("E8").Value = 34469586410
$ws.Range("F8").Value = 930402392
$ws.Range("G8").Value = -0.05828

$ws.Range("D9").Value = 1.001
$ws.Range("E9").Value = 32462239140
$ws.Range("F9").Value = 4824260731
$ws.Range("G9").Value = 0.10131

$ws.Range("D10").Value = 0.207424
$ws.Range("E10").Value = 29797893772
$ws.Range("F10").Value = 2295994232
$ws.Range("G10").Value = -1.51399

$ws.Range("D11").Value = 0.647124
$ws.Range("E11").Value = 22825795620
$ws.Range("F11").Value = 325283919
$ws.Range("G11").Value = -1.52699

$ws.Range("E12").Value = 20299942121
$ws.Range("F12").Value = 307394853
$ws.Range("G12").Value = -0.91729

$ws.Range("D13").Value = 5.17
$ws.Range("E13").Value = 17951052549
$ws.Range("F13").Value = 139606248
$ws.Range("G13").Value = -1.82854

$ws.Range("D14").Value = 0.00003042
$ws.Range("E14").Value = 17925156777
$ws.Range("F14").Value = 513582872
$ws.Range("G14").Value = -0.1436

$ws.Range("D15").Value = 9.58
$ws.Range("E15").Value = 12894315849
$ws.Range("F15").Value = 155410442
$ws.Range("G15").Value = -0.13951

$ws.Range("B16").Value = 'BCH'
$ws.Range("C16").Value = 'Bitcoin Cash'
$ws.Range("D16").Value = 595.89
$ws.Range("E16").Value = 11738234482
$ws.Range("F16").Value = 399114506
$ws.Range("G16").Value = -1.297

$ws.Range("B17").Value = 'LINK'
$ws.Range("C17").Value = 'Chainlink'
$ws.Range("D17").Value = 19.05
$ws.Range("E17").Value = 11183126975
$ws.Range("F17").Value = 362435030
$ws.Range("G17").Value = -0.70952

$ws.Range("D18").Value = 70564
$ws.Range("E18").Value = 10956756486
$ws.Range("F18").Value = 183297402
$ws.Range("G18").Value = 0.44085

$ws.Range("B19").Value = 'TRX'
$ws.Range("C19").Value = 'TRON'
$ws.Range("D19").Value = 0.122627
$ws.Range("E19").Value = 10756666177
$ws.Range("F19").Value = 270646828
$ws.Range("G19").Value = 1.43874

$ws.Range("B20").Value = 'UNI'
$ws.Range("C20").Value = 'Uniswap'
$ws.Range("D20").Value = 13.05
$ws.Range("E20").Value = 9835281890
$ws.Range("F20").Value = 183951923
$ws.Range("G20").Value = 2.52941

$ws.Range("B21").Value = 'MATIC'
$ws.Range("C21").Value = 'Polygon'
$ws.Range("D21").Value = 0.998556
$ws.Range("E21").Value = 9268098271
$ws.Range("F21").Value = 209388759
$ws.Range("G21").Value = -0.20482

$ws.Range("B22").Value = 'ICP'
$ws.Range("C22").Value = 'Internet Computer'
$ws.Range("D22").Value = 17.83
$ws.Range("E22").Value = 8242140330
$ws.Range("F22").Value = 227796927
$ws.Range("G22").Value = -1.24331

$ws.Range("B23").Value = 'LTC'
$ws.Range("C23").Value = 'Litecoin'
$ws.Range("D23").Value = 102.27
$ws.Range("E23").Value = 7612536539
$ws.Range("F23").Value = 588154818
$ws.Range("G23").Value = -1.25229

$ws.Range("B24").Value = 'NEAR'
$ws.Range("C24").Value = 'NEAR Protocol'
$ws.Range("D24").Value = 7.24
$ws.Range("E24").Value = 7508174086
$ws.Range("F24").Value = 237992068
$ws.Range("G24").Value = 1.01461

$ws.Range("B25").Value = 'APT'
$ws.Range("C25").Value = 'Aptos'
$ws.Range("D25").Value = 16.48
$ws.Range("E25").Value = 6543760031
$ws.Range("F25").Value = 143729890
$ws.Range("G25").Value = -3.09986

$ws.Range("D26").Value = 5.97
$ws.Range("E26").Value = 5529214162
$ws.Range("F26").Value = 1824745
$ws.Range("G26").Value = -1.8902

$ws.Range("D27").Value = 3.53
$ws.Range("E27").Value = 5121062940
$ws.Range("F27").Value = 78707648
$ws.Range("G27").Value = -3.59425

$ws.Range("B28").Value = 'FIL'
$ws.Range("C28").Value = 'Filecoin'
$ws.Range("D28").Value = 9.640000000000001
$ws.Range("E28").Value = 5115347533
$ws.Range("F28").Value = 206539209
$ws.Range("G28").Value = -0.71816

$ws.Range("B29").Value = 'ETC'
$ws.Range("C29").Value = 'Ethereum Classic'
$ws.Range("D29").Value = 33.95
$ws.Range("E29").Value = 4963592324
$ws.Range("F29").Value = 166387035
$ws.Range("G29").Value = 0.35056

$ws.Range("B30").Value = 'DAI'
$ws.Range("C30").Value = 'Dai'
$ws.Range("D30").Value = 0.999179
$ws.Range("E30").Value = 4911796100
$ws.Range("F30").Value = 238880024
$ws.Range("G30").Value = -0.02738

$ws.Range("D31").Value = 12.32
$ws.Range("E31").Value = 4805672661
$ws.Range("F31").Value = 121828600
$ws.Range("G31").Value = -3.33235

$ws.Range("B32").Value = 'WIF'
$ws.Range("C32").Value = 'dogwifhat'
$ws.Range("D32").Value = 4.77
$ws.Range("E32").Value = 4765035952
$ws.Range("F32").Value = 624374578
$ws.Range("G32").Value = 5.37514

$ws.Range("B33").Value = 'ARB'
$ws.Range("C33").Value = 'Arbitrum'
$ws.Range("D33").Value = 1.66
$ws.Range("E33").Value = 4388618312
$ws.Range("F33").Value = 236675295
$ws.Range("G33").Value = -1.068

$ws.Range("B34").Value = 'IMX'
$ws.Range("C34").Value = 'Immutable'
$ws.Range("D34").Value = 3.03
$ws.Range("E34").Value = 4313282802
$ws.Range("F34").Value = 59049519
$ws.Range("G34").Value = -3.17973

$ws.Range("B35").Value = 'MNT'
$ws.Range("C35").Value = 'Mantle'
$ws.Range("D35").Value = 1.31
$ws.Range("E35").Value = 4246870182
$ws.Range("F35").Value = 143962000
$ws.Range("G35").Value = -0.27654

$ws.Range("B36").Value = 'RNDR'
$ws.Range("C36").Value = 'Render'
$ws.Range("D36").Value = 10.78
$ws.Range("E36").Value = 4116422908
$ws.Range("F36").Value = 141029275
$ws.Range("G36").Value = -1.83882

$ws.Range("B37").Value = 'CRO'
$ws.Range("C37").Value = 'Cronos'
$ws.Range("D37").Value = 0.152624
$ws.Range("E37").Value = 4061305323
$ws.Range("F37").Value = 14683531
$ws.Range("G37").Value = -0.47143

$ws.Range("B38").Value = 'XLM'
$ws.Range("C38").Value = 'Stellar'
$ws.Range("D38").Value = 0.140521
$ws.Range("E38").Value = 4047702496
$ws.Range("F38").Value = 66304741
$ws.Range("G38").Value = 0.03147

$ws.Range("B39").Value = 'HBAR'
$ws.Range("C39").Value = 'Hedera'
$ws.Range("D39").Value = 0.117449
$ws.Range("E39").Value = 3962373628
$ws.Range("F39").Value = 46647765
$ws.Range("G39").Value = 1.05005

$ws.Range("B40").Value = 'PEPE'
$ws.Range("C40").Value = 'Pepe'
$ws.Range("D40").Value = 0.00000902
$ws.Range("E40").Value = 3799832344
$ws.Range("F40").Value = 1146524984
$ws.Range("G40").Value = 3.94063

$ws.Range("B41").Value = 'OKB'
$ws.Range("C41").Value = 'OKB'
$ws.Range("D41").Value = 63.24
$ws.Range("E41").Value = 3794180368
$ws.Range("F41").Value = 7304948
$ws.Range("G41").Value = -0.41365

$ws.Range("B42").Value = 'OP'
$ws.Range("C42").Value = 'Optimism'
$ws.Range("D42").Value = 3.7
$ws.Range("E42").Value = 3714672850
$ws.Range("F42").Value = 155775414
$ws.Range("G42").Value = -0.70314

$ws.Range("B43").Value = 'GRT'
$ws.Range("C43").Value = 'The Graph'
$ws.Range("D43").Value = 0.39206
$ws.Range("E43").Value = 3705636054
$ws.Range("F43").Value = 74328806
$ws.Range("G43").Value = -1.19136

$ws.Range("B44").Value = 'MKR'
$ws.Range("C44").Value = 'Maker'
$ws.Range("D44").Value = 3894.32
$ws.Range("E44").Value = 3601277038
$ws.Range("F44").Value = 164010407
$ws.Range("G44").Value = 3.54241

$ws.Range("B45").Value = 'TAO'
$ws.Range("C45").Value = 'Bittensor'
$ws.Range("D45").Value = 538.8
$ws.Range("E45").Value = 3496246783
$ws.Range("F45").Value = 58362788
$ws.Range("G45").Value = 9.747629999999999

$ws.Range("B46").Value = 'VET'
$ws.Range("C46").Value = 'VeChain'
$ws.Range("D46").Value = 0.04554033
$ws.Range("E46").Value = 3304719594
$ws.Range("F46").Value = 66237651
$ws.Range("G46").Value = -0.51629

$ws.Range("B47").Value = 'INJ'
$ws.Range("C47").Value = 'Injective'
$ws.Range("D47").Value = 37.01
$ws.Range("E47").Value = 3272486357
$ws.Range("F47").Value = 95670313
$ws.Range("G47").Value = 0.00211

$ws.Range("B48").Value = 'FET'
$ws.Range("C48").Value = 'Fetch.ai'
$ws.Range("D48").Value = 3.12
$ws.Range("E48").Value = 3252841533
$ws.Range("F48").Value = 404909325
$ws.Range("G48").Value = 0.69363

$ws.Range("B49").Value = 'KAS'
$ws.Range("C49").Value = 'Kaspa'
$ws.Range("D49").Value = 0.133446
$ws.Range("E49").Value = 3079428605
$ws.Range("F49").Value = 30092141
$ws.Range("G49").Value = -1.64291

$ws.Range("B50").Value = 'RUNE'
$ws.Range("C50").Value = 'THORChain'
$ws.Range("D50").Value = 8.630000000000001
$ws.Range("E50").Value = 2896186476
$ws.Range("F50").Value = 183133122
$ws.Range("G50").Value = -0.83719

$ws.Range("B51").Value = 'THETA'
$ws.Range("C51").Value = 'Theta Network'
$ws.Range("D51").Value = 2.88
$ws.Range("E51").Value = 2873366488
$ws.Range("F51").Value = 29338379
$ws.Range("G51").Value = 0.66028
